# Atualizei dados bibi e add
# Update faturamento_diario_lojas data: column T (outro meio de pagamento)
# and column AG (total) for rows 2-6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("T2").Value = 9216.52
$ws.Range("AG2").Value = 242762.19

# Row 3
$ws.Range("T3").Value = 3678
$ws.Range("AG3").Value = 88065.3

# Row 4
$ws.Range("T4").Value = 1150
$ws.Range("AG4").Value = 56517.33

# Row 5
$ws.Range("T5").Value = 1503
$ws.Range("AG5").Value = 50507.05

# Row 6
$ws.Range("T6").Value = 15547.52
$ws.Range("AG6").Value = 437851.87
